# Apply cryptos list update (price/volume refresh + rank-47/48 coin swap)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '68.249.54'
$ws.Range('E2').Value = '  +1.12%  '

$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.912.64'
$ws.Range('E3').Value = '  -0.80%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '487.82'
$ws.Range('E5').Value = '  +3.51%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '146.93'
$ws.Range('E6').Value = '  -0.55%  '

$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.623'
$ws.Range('E7').Value = '  -0.34%  '

$ws.Range('E8').Value = '  -0.10%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.729'
$ws.Range('E9').Value = '  -0.45%  '

$ws.Range('E10').Value = '  -1.30%  '

$ws.Range('E11').Value = '  -2.11%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '43.08'
$ws.Range('E12').Value = '  -0.68%  '

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '10.91'
$ws.Range('E13').Value = '  +4.87%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '4.538.61'
$ws.Range('E14').Value = '  -0.70%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '3.912.87'
$ws.Range('E15').Value = '  -0.87%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '14.24'
$ws.Range('E16').Value = '  -5.78%  '

$ws.Range('E17').Value = '  -1.18%  '

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '19.93'
$ws.Range('E18').Value = '  -0.03%  '

$ws.Range('E19').Value = '  -1.81%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '68.359.31'
$ws.Range('E20').Value = '  +1.04%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '434.02'
$ws.Range('E21').Value = '  -0.29%  '

$ws.Range('E22').Value = '  +4.45%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '14.86'
$ws.Range('E23').Value = '  +2.76%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '87.89'
$ws.Range('E24').Value = '  +0.45%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '11.39'
$ws.Range('E25').Value = '  +14.59%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '11.25'
$ws.Range('E26').Value = '  +10.66%  '

$ws.Range('E27').Value = '  +0.17%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '38.10'
$ws.Range('E28').Value = '  -1.28%  '

$ws.Range('E29').Value = '  +1.14%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '722.74'
$ws.Range('E30').Value = '  -0.07%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '13.77'
$ws.Range('E31').Value = '  +2.55%  '

$ws.Range('E32').Value = '  -1.89%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '2.93'
$ws.Range('E33').Value = '  +4.16%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '6.27'
$ws.Range('E34').Value = '  +16.97%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '41.66'
$ws.Range('E35').Value = '  -1.62%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.0₃0875'
$ws.Range('E36').Value = '  +4.64%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '60.31'
$ws.Range('E37').Value = '  +4.08%  '

$ws.Range('E38').Value = '  +19.60%  '

$ws.Range('E39').Value = '  -2.40%  '

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.998'
$ws.Range('E40').Value = '  -0.15%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '2.94'
$ws.Range('E41').Value = '  +15.33%  '

$ws.Range('E42').Value = '  +1.20%  '

$ws.Range('E43').Value = '  +3.40%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '2.91'
$ws.Range('E44').Value = '  +3.25%  '

$ws.Range('E45').Value = '  -0.62%  '

$ws.Range('E46').Value = '  +0.07%  '

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '3.33'
$ws.Range('E47').Value = '  +1.41%  '

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '3.41'
$ws.Range('E48').Value = '  -2.50%  '

$ws.Range('B49').Value = 'BabyDogeCoin'
$ws.Range('C49').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.0₆0350'
$ws.Range('E49').Value = '  +37.23%  '

$ws.Range('B50').Value = 'ARBITRUM'
$ws.Range('C50').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '2.14'
$ws.Range('E50').Value = '  -3.50%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '144.63'
$ws.Range('E51').Value = '  -1.90%  '
